$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format on the changed Price/Volume cells so that
# numeric-looking strings (e.g. "291.64", "0.06770", "1.571.10") are written back
# verbatim as text rather than being coerced into floating point numbers.
$cells = @{
    'D2' = '22.474.86'
    'E2' = '  +0.48%  '
    'D3' = '1.572.33'
    'E3' = '  +0.31%  '
    'E5' = '  -0.06%  '
    'D6' = '291.64'
    'E6' = '  +0.14%  '
    'D7' = '0.3712'
    'E7' = '  -1.34%  '
    'D8' = '49.99'
    'E8' = '  +0.92%  '
    'D9' = '0.3389'
    'E9' = '  -0.53%  '
    'D10' = '1.143'
    'E10' = '  +0.09%  '
    'D11' = '0.07547'
    'E11' = '  -0.73%  '
    'E12' = '  -0.08%  '
    'D13' = '21.29'
    'E13' = '  +1.11%  '
    'D14' = '6.038'
    'E14' = '  +0.83%  '
    'D15' = '6.959'
    'E15' = '  +0.06%  '
    'D16' = '1.571.10'
    'E16' = '  +0.18%  '
    'E17' = '  -0.81%  '
    'D18' = '90.66'
    'E18' = '  +0.75%  '
    'D19' = '0.06770'
    'D20' = '1.002'
    'E20' = '  -0.06%  '
    'D21' = '6.288'
    'E21' = '  +1.57%  '
    'D22' = '16.38'
    'E22' = '  -0.96%  '
    'E23' = '  +2.10%  '
    'D24' = '22.464.89'
    'E24' = '  +0.48%  '
    'D25' = '2.355'
    'E25' = '  -2.08%  '
    'D26' = '2.621'
    'E26' = '  -3.47%  '
    'E27' = '  -0.39%  '
    'D28' = '149.47'
    'E28' = '  +1.38%  '
    'D29' = '5.064'
    'E29' = '  +0.65%  '
    'D30' = '125.08'
    'E30' = '  -1.10%  '
    'D31' = '1.746.48'
    'E31' = '  +0.12%  '
    'D32' = '1.083'
    'E32' = '  +8.01%  '
    'D33' = '6.202'
    'E33' = '  +1.67%  '
    'D34' = '2.009'
    'E34' = '  -0.45%  '
    'D35' = '9.796'
    'E35' = '  -2.87%  '
    'D36' = '0.08368'
    'E36' = '  -1.59%  '
    'D37' = '0.02482'
    'E37' = '  -1.37%  '
    'D38' = '0.2301'
    'E38' = '  -0.07%  '
    'D39' = '1.341'
    'E39' = '  -6.48%  '
    'D40' = '0.06545'
    'E40' = '  +0.78%  '
    'D41' = '5.453'
    'E41' = '  +0.99%  '
    'D42' = '11.38'
    'E42' = '  +0.16%  '
    'D43' = '0.6223'
    'E43' = '  -1.63%  '
    'E44' = '  -0.06%  '
    'D45' = '14.05'
    'E45' = '  +0.36%  '
    'D46' = '3.813'
    'E46' = '  +0.30%  '
    'D47' = '0.5853'
    'E47' = '  -1.41%  '
    'D48' = '129.51'
    'E48' = '  +4.20%  '
    'D49' = '2.071'
    'E49' = '  -0.31%  '
    'D50' = '1.215'
    'E50' = '  -5.17%  '
    'D51' = '0.07336'
    'E51' = '  +0.34%  '
}

foreach ($addr in $cells.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $cells[$addr]
}
